$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (73) with the next day's gold price data, following the
# same pattern as the existing rows (date in column A, price text in column B).
$ws.Range("A73").Value2 = "28-11-2025"
$ws.Range("B73").Value2 = "The price of gold in India today is ₹12,846 per gram for 24 karat gold, ₹11,775 per gram for 22 karat gold and ₹9,634 per gram for 18 karat gold (also called 999 gold)."
